$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24:89 down to 25:90.
$ws.Rows.Item(24).EntireRow.Insert()

# Populate the newly inserted row 24 with the new price-report record.
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = 45099
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 100112013
$ws.Range("G24").Value = "Alcachofa"
$ws.Range("H24").Value = "Argentina(o)"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 130
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 17000
$ws.Range("M24").Value = 16615
$ws.Range("N24").Value = "$/caja 50 unidades"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 332
$ws.Range("Q24").Value = 50
$ws.Range("R24").Value = "Hortaliza"
